# Rename the "_old"/"_new" column header suffixes to the respective
# format-version identifiers ("_FV2310" / "_FV2404") used as input file
# names, then turn the header range into a proper Excel Table and freeze
# the header row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns A:J describe the "old" input file, columns L:U the "new" one
# (column K just holds the literal "diff" marker and stays untouched).
$lastCol = 21
for ($col = 1; $col -le $lastCol; $col++) {
    $cell = $ws.Cells.Item(1, $col)
    $header = $cell.Value()
    if ($header -eq $null) { continue }
    $header = [string]$header
    if ($header.EndsWith("_old")) {
        $cell.Value = $header.Substring(0, $header.Length - 4) + "_FV2310"
    } elseif ($header.EndsWith("_new")) {
        $cell.Value = $header.Substring(0, $header.Length - 4) + "_FV2404"
    }
}

# Turn the data range into a proper Excel Table ("Table1") with an
# autofilter, matching the exported AHB-diff layout.
$headerRange = $ws.Range("A1:U64")
$tbl = $ws.ListObjects.Add([Microsoft.Office.Interop.Excel.XlListObjectSourceType]::xlSrcRange, $headerRange, $null, [Microsoft.Office.Interop.Excel.XlYesNoGuess]::xlYes)
$tbl.Name = "Table1"

# Freeze the header row (split after row 1).
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
